# Updated cryptos list with GitHub Actions — refresh price / volume(1h) figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a "Price" (column D) cell, forcing Text format so that values
# which look numeric (e.g. "20.70", "4.500", "0.000008945") keep their exact
# literal representation instead of being normalized/re-formatted by Excel.
function Set-PriceCell($row, $text) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

function Set-VolumeCell($row, $text) {
    $ws.Range("E$row").Value = $text
}

# Row 2 - Bitcoin
Set-PriceCell 2 "27.677.45"
Set-VolumeCell 2 "  +0.18%  "

# Row 3 - Ethereum
Set-PriceCell 3 "1.848.22"
Set-VolumeCell 3 "  +0.17%  "

# Row 4 - TetherUSD
Set-PriceCell 4 "1.006"
Set-VolumeCell 4 "  +0.59%  "

# Row 5 - BNB
Set-PriceCell 5 "312.61"
Set-VolumeCell 5 "  -0.97%  "

# Row 6 - USDC (price unchanged, only volume)
Set-VolumeCell 6 "  +0.47%  "

# Row 7 - XRP
Set-PriceCell 7 "0.4271"
Set-VolumeCell 7 "  +0.59%  "

# Row 8 - Cardano
Set-PriceCell 8 "0.3571"
Set-VolumeCell 8 "  -1.98%  "

# Row 9 - Dogecoin
Set-PriceCell 9 "0.07296"
Set-VolumeCell 9 "  +0.20%  "

# Row 10 - Polygon
Set-PriceCell 10 "0.8695"
Set-VolumeCell 10 "  -2.37%  "

# Row 11 - Solana
Set-PriceCell 11 "20.70"
Set-VolumeCell 11 "  -0.04%  "

# Row 12 - WrappedEther
Set-PriceCell 12 "1.847.05"
Set-VolumeCell 12 "  +1.40%  "

# Row 13 - Chainlink
Set-PriceCell 13 "6.542"
Set-VolumeCell 13 "  -0.39%  "

# Row 14 - Polkadot
Set-PriceCell 14 "5.331"
Set-VolumeCell 14 "  -0.47%  "

# Row 15 - TRON
Set-PriceCell 15 "0.06988"
Set-VolumeCell 15 "  +1.73%  "

# Row 16 - BinanceUSD (price unchanged, only volume)
Set-VolumeCell 16 "  +0.54%  "

# Row 17 - Litecoin
Set-PriceCell 17 "79.62"
Set-VolumeCell 17 "  +0.59%  "

# Row 18 - ShibaInu
Set-PriceCell 18 "0.000008945"
Set-VolumeCell 18 "  +0.73%  "

# Row 19 - Dai
Set-PriceCell 19 "1.004"
Set-VolumeCell 19 "  +0.44%  "

# Row 20 - Avalanche
Set-PriceCell 20 "15.28"
Set-VolumeCell 20 "  -1.32%  "

# Row 21 - WrappedBTC
Set-PriceCell 21 "27.731.05"
Set-VolumeCell 21 "  +0.40%  "

# Row 22 - Uniswap
Set-PriceCell 22 "4.991"
Set-VolumeCell 22 "  -0.06%  "

# Row 23 - Cosmos
Set-PriceCell 23 "10.36"
Set-VolumeCell 23 "  -2.19%  "

# Row 24 - WrappedliquidstakedEther2.0
Set-PriceCell 24 "2.125.07"
Set-VolumeCell 24 "  +4.50%  "

# Row 25 - Toncoin
Set-PriceCell 25 "1.985"
Set-VolumeCell 25 "  +2.07%  "

# Row 26 - Monero (volume unchanged, only price)
Set-PriceCell 26 "155.53"

# Row 27 - EthereumClassic (price unchanged, only volume)
Set-VolumeCell 27 "  -2.79%  "

# Row 28 - BitcoinCash
Set-PriceCell 28 "120.44"
Set-VolumeCell 28 "  -1.42%  "

# Row 29 - InternetComputer(DFINITY)
Set-PriceCell 29 "5.262"
Set-VolumeCell 29 "  -0.43%  "

# Row 30 - LidoDAOToken
Set-PriceCell 30 "1.865"
Set-VolumeCell 30 "  -0.52%  "

# Row 31 - Stellar
Set-PriceCell 31 "0.08911"
Set-VolumeCell 31 "  -0.22%  "

# Row 32 - ImmutableX
Set-PriceCell 32 "0.7582"
Set-VolumeCell 32 "  -1.88%  "

# Row 33 - HuobiToken
Set-PriceCell 33 "2.971"
Set-VolumeCell 33 "  +1.82%  "

# Row 34 - Filecoin
Set-PriceCell 34 "4.500"
Set-VolumeCell 34 "  -1.93%  "

# Row 35 - ARBITRUM
Set-PriceCell 35 "1.124"
Set-VolumeCell 35 "  +2.29%  "

# Row 36 - Frax (price unchanged, only volume)
Set-VolumeCell 36 "  +0.47%  "

# Row 37 - Hedera
Set-PriceCell 37 "0.05432"
Set-VolumeCell 37 "  +0.95%  "

# Row 38 - TrustWalletToken
Set-PriceCell 38 "1.103"
Set-VolumeCell 38 "  +0.32%  "

# Row 39 - VeChain
Set-PriceCell 39 "0.01928"
Set-VolumeCell 39 "  -0.27%  "

# Row 40 - MXToken
Set-PriceCell 40 "2.835"
Set-VolumeCell 40 "  +1.33%  "

# Row 41 - Algorand
Set-PriceCell 41 "0.1664"
Set-VolumeCell 41 "  +0.55%  "

# Row 42 - TheSandbox
Set-PriceCell 42 "0.5067"
Set-VolumeCell 42 "  -0.79%  "

# Row 43 - FraxShare
Set-PriceCell 43 "6.595"
Set-VolumeCell 43 "  -4.30%  "

# Row 44 - Aptos
Set-PriceCell 44 "8.384"
Set-VolumeCell 44 "  +0.89%  "

# Row 45 - was Quant, now Cronos (coin/link/price/volume all swap)
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-PriceCell 45 "0.06546"
Set-VolumeCell 45 "  -0.59%  "

# Row 46 - was Cronos, now Quant (coin/link/price/volume all swap)
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-PriceCell 46 "106.38"
Set-VolumeCell 46 "  +1.35%  "

# Row 47 - EnergySwap
Set-PriceCell 47 "10.37"
Set-VolumeCell 47 "  -0.10%  "

# Row 48 - Decentraland
Set-PriceCell 48 "0.4653"
Set-VolumeCell 48 "  -1.80%  "

# Row 49 - PaxDollar (price unchanged, only volume)
Set-VolumeCell 49 "  +0.48%  "

# Row 50 - NEARProtocol
Set-PriceCell 50 "1.628"
Set-VolumeCell 50 "  -0.59%  "

# Row 51 - Aave
Set-PriceCell 51 "64.42"
Set-VolumeCell 51 "  -0.19%  "
